# Update "想去人数" (F column) figures on the "展览", "演出" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 280
$ws1.Range("F4").Value = 7730
$ws1.Range("F5").Value = 5623
$ws1.Range("F6").Value = 466
$ws1.Range("F7").Value = 74
$ws1.Range("F10").Value = 252
$ws1.Range("F11").Value = 221
$ws1.Range("F12").Value = 55

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 84

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 280
$ws4.Range("F4").Value = 7730
$ws4.Range("F5").Value = 5623
$ws4.Range("F6").Value = 466
$ws4.Range("F7").Value = 74
$ws4.Range("F10").Value = 252
$ws4.Range("F11").Value = 84
$ws4.Range("F13").Value = 221
$ws4.Range("F14").Value = 55
